# Fix bug: not removing ID and TIME in features. Rerun training and model assessment.
#
# The training pipeline was leaking the row identifier (ID_APPLICATION) into the
# feature set used to fit LogisticRegression. After excluding ID/TIME-derived
# columns from training, the model was retrained and this feature-importance
# report was regenerated: the ID_APPLICATION row is gone and every remaining
# feature now carries its refreshed abs_importance / importance (and, since the
# sheet is sorted by abs_importance, rows shift position accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the leaked ID feature row; everything below shifts up by one.
$idRow = 0
$idCell = $ws.Columns.Item(1).Find("ID_APPLICATION")
if ($idCell -ne $null) {
    $idRow = $idCell.Row
} else {
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 1).Value2 -eq "ID_APPLICATION") {
            $idRow = $r
            break
        }
    }
}
if ($idRow -gt 0) {
    $ws.Rows.Item($idRow).Delete()
}

# Rewrite the refreshed, re-sorted feature report (rows 2..165) in one shot.
$data = New-Object 'object[,]' 164,3
$data[0,0] = "NUMERICAL_4"
$data[0,1] = 1.972867202374893
$data[0,2] = -1.972867202374893
$data[1,0] = "NUMERICAL_4_std_dev_last_30_days"
$data[1,1] = 1.971266631690917
$data[1,2] = 1.971266631690917
$data[2,0] = "NUMERICAL_20_std_dev_last_30_days"
$data[2,1] = 1.878460343696272
$data[2,2] = -1.878460343696272
$data[3,0] = "NUMERICAL_20"
$data[3,1] = 1.876466050243653
$data[3,2] = 1.876466050243653
$data[4,0] = "NUMERICAL_18"
$data[4,1] = 1.711722455031699
$data[4,2] = 1.711722455031699
$data[5,0] = "NUMERICAL_18_std_dev_last_30_days"
$data[5,1] = 1.705760223294286
$data[5,2] = -1.705760223294286
$data[6,0] = "NUMERICAL_7"
$data[6,1] = 1.666467294751642
$data[6,2] = -1.666467294751642
$data[7,0] = "NUMERICAL_7_std_dev_last_30_days"
$data[7,1] = 1.659492596044235
$data[7,2] = 1.659492596044235
$data[8,0] = "NUMERICAL_8_std_dev_last_30_days"
$data[8,1] = 1.347239778896966
$data[8,2] = 1.347239778896966
$data[9,0] = "NUMERICAL_8"
$data[9,1] = 1.347000357052889
$data[9,2] = -1.347000357052889
$data[10,0] = "NUMERICAL_32_std_dev_last_30_days"
$data[10,1] = 1.254475242266157
$data[10,2] = -1.254475242266157
$data[11,0] = "NUMERICAL_32"
$data[11,1] = 1.252158478164886
$data[11,2] = 1.252158478164886
$data[12,0] = "NUMERICAL_11"
$data[12,1] = 1.250186222938902
$data[12,2] = 1.250186222938902
$data[13,0] = "NUMERICAL_11_std_dev_last_30_days"
$data[13,1] = 1.246855862823902
$data[13,2] = -1.246855862823902
$data[14,0] = "NUMERICAL_28"
$data[14,1] = 1.160202461512431
$data[14,2] = -1.160202461512431
$data[15,0] = "NUMERICAL_28_std_dev_last_30_days"
$data[15,1] = 1.149869891247123
$data[15,2] = 1.149869891247123
$data[16,0] = "NUMERICAL_6"
$data[16,1] = 1.114194718010076
$data[16,2] = -1.114194718010076
$data[17,0] = "NUMERICAL_6_std_dev_last_30_days"
$data[17,1] = 1.112403182527738
$data[17,2] = 1.112403182527738
$data[18,0] = "NUMERICAL_40_std_dev_last_30_days"
$data[18,1] = 1.069885245401764
$data[18,2] = 1.069885245401764
$data[19,0] = "NUMERICAL_39_std_dev_last_30_days"
$data[19,1] = 1.02160183504064
$data[19,2] = 1.02160183504064
$data[20,0] = "NUMERICAL_39"
$data[20,1] = 1.008630107904062
$data[20,2] = -1.008630107904062
$data[21,0] = "NUMERICAL_16_std_dev_last_30_days"
$data[21,1] = 1.003527321239619
$data[21,2] = 1.003527321239619
$data[22,0] = "NUMERICAL_16"
$data[22,1] = 0.9912640850371018
$data[22,2] = -0.9912640850371018
$data[23,0] = "NUMERICAL_41_std_dev_last_30_days"
$data[23,1] = 0.9857034989990721
$data[23,2] = 0.9857034989990721
$data[24,0] = "NUMERICAL_29"
$data[24,1] = 0.9490226336024571
$data[24,2] = -0.9490226336024571
$data[25,0] = "NUMERICAL_29_std_dev_last_30_days"
$data[25,1] = 0.9385550747022456
$data[25,2] = 0.9385550747022456
$data[26,0] = "NUMERICAL_36"
$data[26,1] = 0.9043446853245833
$data[26,2] = 0.9043446853245833
$data[27,0] = "NUMERICAL_36_std_dev_last_30_days"
$data[27,1] = 0.8981198747201898
$data[27,2] = -0.8981198747201898
$data[28,0] = "NUMERICAL_9"
$data[28,1] = 0.8381534343579177
$data[28,2] = -0.8381534343579177
$data[29,0] = "NUMERICAL_9_std_dev_last_30_days"
$data[29,1] = 0.8241306187908757
$data[29,2] = 0.8241306187908757
$data[30,0] = "NUMERICAL_38"
$data[30,1] = 0.7539077161022443
$data[30,2] = 0.7539077161022443
$data[31,0] = "NUMERICAL_38_std_dev_last_30_days"
$data[31,1] = 0.752516268852252
$data[31,2] = -0.752516268852252
$data[32,0] = "NUMERICAL_17"
$data[32,1] = 0.735110797969498
$data[32,2] = -0.735110797969498
$data[33,0] = "NUMERICAL_17_std_dev_last_30_days"
$data[33,1] = 0.7292560578262479
$data[33,2] = 0.7292560578262479
$data[34,0] = "NUMERICAL_25"
$data[34,1] = 0.7178408600855275
$data[34,2] = -0.7178408600855275
$data[35,0] = "NUMERICAL_25_std_dev_last_30_days"
$data[35,1] = 0.7130565062589005
$data[35,2] = 0.7130565062589005
$data[36,0] = "NUMERICAL_40"
$data[36,1] = 0.6819595679567425
$data[36,2] = -0.6819595679567425
$data[37,0] = "NUMERICAL_5_std_dev_last_30_days"
$data[37,1] = 0.6697698805372652
$data[37,2] = 0.6697698805372652
$data[38,0] = "NUMERICAL_5"
$data[38,1] = 0.6639868535600477
$data[38,2] = -0.6639868535600477
$data[39,0] = "CATEGORICAL_7_value_C"
$data[39,1] = 0.6278098002294416
$data[39,2] = -0.6278098002294416
$data[40,0] = "NUMERICAL_37_std_dev_last_30_days"
$data[40,1] = 0.6231706036827964
$data[40,2] = -0.6231706036827964
$data[41,0] = "NUMERICAL_31"
$data[41,1] = 0.6216686731586589
$data[41,2] = 0.6216686731586589
$data[42,0] = "NUMERICAL_31_std_dev_last_30_days"
$data[42,1] = 0.6160292075661719
$data[42,2] = -0.6160292075661719
$data[43,0] = "NUMERICAL_37"
$data[43,1] = 0.6143352796792046
$data[43,2] = 0.6143352796792046
$data[44,0] = "CATEGORICAL_7_value_A"
$data[44,1] = 0.6014448571174996
$data[44,2] = 0.6014448571174996
$data[45,0] = "NUMERICAL_3_std_dev_last_30_days"
$data[45,1] = 0.5917138277664491
$data[45,2] = 0.5917138277664491
$data[46,0] = "NUMERICAL_30"
$data[46,1] = 0.5844352341549742
$data[46,2] = 0.5844352341549742
$data[47,0] = "NUMERICAL_3"
$data[47,1] = 0.5824969255029668
$data[47,2] = -0.5824969255029668
$data[48,0] = "CATEGORICAL_9_value_AA"
$data[48,1] = 0.5618641816389968
$data[48,2] = 0.5618641816389968
$data[49,0] = "CATEGORICAL_9_value_BB"
$data[49,1] = 0.5618641816389968
$data[49,2] = -0.5618641816389968
$data[50,0] = "NUMERICAL_24_std_dev_last_30_days"
$data[50,1] = 0.5092927852142224
$data[50,2] = 0.5092927852142224
$data[51,0] = "NUMERICAL_24"
$data[51,1] = 0.5003962368892488
$data[51,2] = -0.5003962368892488
$data[52,0] = "CATEGORICAL_1_value_DVSOM"
$data[52,1] = 0.4538522834820553
$data[52,2] = -0.4538522834820553
$data[53,0] = "NUMERICAL_13_std_dev_last_30_days"
$data[53,1] = 0.4527351436223263
$data[53,2] = -0.4527351436223263
$data[54,0] = "NUMERICAL_21"
$data[54,1] = 0.4459628754566192
$data[54,2] = -0.4459628754566192
$data[55,0] = "NUMERICAL_13"
$data[55,1] = 0.4394745195813375
$data[55,2] = 0.4394745195813375
$data[56,0] = "NUMERICAL_21_std_dev_last_30_days"
$data[56,1] = 0.4340748167165208
$data[56,2] = 0.4340748167165208
$data[57,0] = "NUMERICAL_22_std_dev_last_30_days"
$data[57,1] = 0.4290963839474641
$data[57,2] = -0.4290963839474641
$data[58,0] = "NUMERICAL_22"
$data[58,1] = 0.4273290704760714
$data[58,2] = 0.4273290704760714
$data[59,0] = "NUMERICAL_34"
$data[59,1] = 0.3972231697560624
$data[59,2] = 0.3972231697560624
$data[60,0] = "NUMERICAL_19_std_dev_last_30_days"
$data[60,1] = 0.3963320849269483
$data[60,2] = 0.3963320849269483
$data[61,0] = "NUMERICAL_34_std_dev_last_30_days"
$data[61,1] = 0.3888820508079136
$data[61,2] = -0.3888820508079136
$data[62,0] = "NUMERICAL_19"
$data[62,1] = 0.3884650075334233
$data[62,2] = -0.3884650075334233
$data[63,0] = "NUMERICAL_15"
$data[63,1] = 0.3799393940476218
$data[63,2] = -0.3799393940476218
$data[64,0] = "NUMERICAL_15_std_dev_last_30_days"
$data[64,1] = 0.3648779972671397
$data[64,2] = 0.3648779972671397
$data[65,0] = "NUMERICAL_35_std_dev_last_30_days"
$data[65,1] = 0.3625654406160402
$data[65,2] = -0.3625654406160402
$data[66,0] = "NUMERICAL_35"
$data[66,1] = 0.3592474067357506
$data[66,2] = 0.3592474067357506
$data[67,0] = "NUMERICAL_41"
$data[67,1] = 0.3305932266251619
$data[67,2] = -0.3305932266251619
$data[68,0] = "CATEGORICAL_1_value_OMFDE"
$data[68,1] = 0.2742304582254372
$data[68,2] = 0.2742304582254372
$data[69,0] = "NUMERICAL_14_std_dev_last_30_days"
$data[69,1] = 0.2692701445749445
$data[69,2] = -0.2692701445749445
$data[70,0] = "NUMERICAL_14"
$data[70,1] = 0.2597235079340913
$data[70,2] = 0.2597235079340913
$data[71,0] = "NUMERICAL_23"
$data[71,1] = 0.2456697249782782
$data[71,2] = 0.2456697249782782
$data[72,0] = "NUMERICAL_23_std_dev_last_30_days"
$data[72,1] = 0.2335582782641035
$data[72,2] = -0.2335582782641035
$data[73,0] = "NUMERICAL_12_std_dev_last_30_days"
$data[73,1] = 0.2224650746959
$data[73,2] = -0.2224650746959
$data[74,0] = "NUMERICAL_12"
$data[74,1] = 0.2209614849489828
$data[74,2] = 0.2209614849489828
$data[75,0] = "CATEGORICAL_4_value_GH"
$data[75,1] = 0.2175294261958693
$data[75,2] = -0.2175294261958693
$data[76,0] = "CATEGORICAL_2_value_AA"
$data[76,1] = 0.2150177212267913
$data[76,2] = -0.2150177212267913
$data[77,0] = "CATEGORICAL_5_value_TR"
$data[77,1] = 0.2117630424375003
$data[77,2] = -0.2117630424375003
$data[78,0] = "CATEGORICAL_3_value_AS"
$data[78,1] = 0.2096048283728223
$data[78,2] = -0.2096048283728223
$data[79,0] = "NUMERICAL_1"
$data[79,1] = 0.1865699568957469
$data[79,2] = -0.1865699568957469
$data[80,0] = "NUMERICAL_33"
$data[80,1] = 0.1862217390326682
$data[80,2] = 0.1862217390326682
$data[81,0] = "NUMERICAL_33_std_dev_last_30_days"
$data[81,1] = 0.1804389322857882
$data[81,2] = -0.1804389322857882
$data[82,0] = "CATEGORICAL_1_value_FCOHQ"
$data[82,1] = 0.1791085776735963
$data[82,2] = 0.1791085776735963
$data[83,0] = "NUMERICAL_1_std_dev_last_30_days"
$data[83,1] = 0.1754351193164419
$data[83,2] = 0.1754351193164419
$data[84,0] = "NUMERICAL_0"
$data[84,1] = 0.1633573235391143
$data[84,2] = 0.1633573235391143
$data[85,0] = "NUMERICAL_0_std_dev_last_30_days"
$data[85,1] = 0.1534895112681541
$data[85,2] = -0.1534895112681541
$data[86,0] = "CATEGORICAL_4_value_JK"
$data[86,1] = 0.1517322805563497
$data[86,2] = 0.1517322805563497
$data[87,0] = "CATEGORICAL_2_value_AB"
$data[87,1] = 0.148528432444057
$data[87,2] = 0.148528432444057
$data[88,0] = "CATEGORICAL_5_value_WE"
$data[88,1] = 0.1456144342492815
$data[88,2] = 0.1456144342492815
$data[89,0] = "CATEGORICAL_3_value_DF"
$data[89,1] = 0.1430319639467789
$data[89,2] = 0.1430319639467789
$data[90,0] = "MONTH_APPLICATION_value_4"
$data[90,1] = 0.1429576498229529
$data[90,2] = 0.1429576498229529
$data[91,0] = "NUMERICAL_30_std_dev_last_30_days"
$data[91,1] = 0.115738025978838
$data[91,2] = -0.115738025978838
$data[92,0] = "NUMERICAL_27_std_dev_last_30_days"
$data[92,1] = 0.1116853756200428
$data[92,2] = 0.1116853756200428
$data[93,0] = "NUMERICAL_27"
$data[93,1] = 0.1019126342553063
$data[93,2] = -0.1019126342553063
$data[94,0] = "CATEGORICAL_5_value_unkown"
$data[94,1] = 0.06196539284892139
$data[94,2] = 0.06196539284892139
$data[95,0] = "CATEGORICAL_2_value_unkown"
$data[95,1] = 0.06196539284892139
$data[95,2] = 0.06196539284892139
$data[96,0] = "CATEGORICAL_4_value_unkown"
$data[96,1] = 0.06196539284892139
$data[96,2] = 0.06196539284892139
$data[97,0] = "CATEGORICAL_3_value_unkown"
$data[97,1] = 0.06196539284892139
$data[97,2] = 0.06196539284892139
$data[98,0] = "CATEGORICAL_7_value_B"
$data[98,1] = 0.05777140100133461
$data[98,2] = -0.05777140100133461
$data[99,0] = "NUMERICAL_26_std_dev_last_30_days"
$data[99,1] = 0.03847486441576235
$data[99,2] = -0.03847486441576235
$data[100,0] = "MONTH_APPLICATION_value_10"
$data[100,1] = 0.03344204937114987
$data[100,2] = -0.03344204937114987
$data[101,0] = "MONTH_APPLICATION_value_11"
$data[101,1] = 0.03246069841469804
$data[101,2] = -0.03246069841469804
$data[102,0] = "CATEGORICAL_8_value_BB"
$data[102,1] = 0.03195790575928951
$data[102,2] = -0.03195790575928951
$data[103,0] = "CATEGORICAL_8_value_AA"
$data[103,1] = 0.03195790575923354
$data[103,2] = 0.03195790575923354
$data[104,0] = "MONTH_APPLICATION_value_12"
$data[104,1] = 0.02883134919428099
$data[104,2] = -0.02883134919428099
$data[105,0] = "NUMERICAL_26"
$data[105,1] = 0.02765454943582196
$data[105,2] = 0.02765454943582196
$data[106,0] = "CATEGORICAL_6_value_A"
$data[106,1] = 0.02241054816895291
$data[106,2] = 0.02241054816895291
$data[107,0] = "DAY_APPLICATION_value_26"
$data[107,1] = 0.02086763077713635
$data[107,2] = 0.02086763077713635
$data[108,0] = "DAY_APPLICATION_value_20"
$data[108,1] = 0.01882588606905898
$data[108,2] = 0.01882588606905898
$data[109,0] = "MONTH_APPLICATION_value_8"
$data[109,1] = 0.01851206534126017
$data[109,2] = -0.01851206534126017
$data[110,0] = "DAY_APPLICATION_value_7"
$data[110,1] = 0.01674489274354554
$data[110,2] = -0.01674489274354554
$data[111,0] = "MONTH_APPLICATION_value_2"
$data[111,1] = 0.01605152028203273
$data[111,2] = -0.01605152028203273
$data[112,0] = "DAY_APPLICATION_value_18"
$data[112,1] = 0.01488588199829113
$data[112,2] = -0.01488588199829113
$data[113,0] = "DAY_APPLICATION_value_15"
$data[113,1] = 0.01479444755104834
$data[113,2] = -0.01479444755104834
$data[114,0] = "CATEGORICAL_0_value_SUPSY"
$data[114,1] = 0.01438254591533904
$data[114,2] = 0.01438254591533904
$data[115,0] = "CATEGORICAL_6_value_C"
$data[115,1] = 0.01423029093649983
$data[115,2] = -0.01423029093649983
$data[116,0] = "DAY_APPLICATION_value_27"
$data[116,1] = 0.01389489668370205
$data[116,2] = 0.01389489668370205
$data[117,0] = "DAY_APPLICATION_value_10"
$data[117,1] = 0.01340163703762189
$data[117,2] = -0.01340163703762189
$data[118,0] = "DAY_APPLICATION_value_28"
$data[118,1] = 0.01312878728968315
$data[118,2] = 0.01312878728968315
$data[119,0] = "DAY_APPLICATION_value_22"
$data[119,1] = 0.01263617083263885
$data[119,2] = 0.01263617083263885
$data[120,0] = "DOW_APPLICATION_value_6"
$data[120,1] = 0.01144724018516868
$data[120,2] = 0.01144724018516868
$data[121,0] = "MONTH_APPLICATION_value_7"
$data[121,1] = 0.01134447957979535
$data[121,2] = -0.01134447957979535
$data[122,0] = "DAY_APPLICATION_value_19"
$data[122,1] = 0.01059549047665766
$data[122,2] = 0.01059549047665766
$data[123,0] = "DAY_APPLICATION_value_23"
$data[123,1] = 0.01032194048784713
$data[123,2] = 0.01032194048784713
$data[124,0] = "CATEGORICAL_6_value_B"
$data[124,1] = 0.01024098103359272
$data[124,2] = -0.01024098103359272
$data[125,0] = "DAY_APPLICATION_value_3"
$data[125,1] = 0.009673989907537033
$data[125,2] = -0.009673989907537033
$data[126,0] = "DAY_APPLICATION_value_17"
$data[126,1] = 0.009551418273577939
$data[126,2] = -0.009551418273577939
$data[127,0] = "CATEGORICAL_0_value_JJUFY"
$data[127,1] = 0.008708396945942956
$data[127,2] = -0.008708396945942956
$data[128,0] = "DAY_APPLICATION_value_30"
$data[128,1] = 0.00856881043097388
$data[128,2] = 0.00856881043097388
$data[129,0] = "HOUR_APPLICATION"
$data[129,1] = 0.008543025478490614
$data[129,2] = 0.008543025478490614
$data[130,0] = "NUMERICAL_2"
$data[130,1] = 0.007885009333720404
$data[130,2] = 0.007885009333720404
$data[131,0] = "DAY_APPLICATION_value_12"
$data[131,1] = 0.007524565653084275
$data[131,2] = -0.007524565653084275
$data[132,0] = "DAY_APPLICATION_value_5"
$data[132,1] = 0.007488313344923793
$data[132,2] = -0.007488313344923793
$data[133,0] = "DAY_APPLICATION_value_16"
$data[133,1] = 0.007308927535851692
$data[133,2] = -0.007308927535851692
$data[134,0] = "DOW_APPLICATION_value_5"
$data[134,1] = 0.007210241773926858
$data[134,2] = -0.007210241773926858
$data[135,0] = "NUMERICAL_2_std_dev_last_30_days"
$data[135,1] = 0.00715661846311132
$data[135,2] = -0.00715661846311132
$data[136,0] = "DAY_APPLICATION_value_6"
$data[136,1] = 0.006392529581572086
$data[136,2] = -0.006392529581572086
$data[137,0] = "DAY_APPLICATION_value_25"
$data[137,1] = 0.006010291200453858
$data[137,2] = 0.006010291200453858
$data[138,0] = "DAY_APPLICATION_value_2"
$data[138,1] = 0.005932472227481251
$data[138,2] = -0.005932472227481251
$data[139,0] = "DAY_APPLICATION_value_29"
$data[139,1] = 0.005764044699607888
$data[139,2] = 0.005764044699607888
$data[140,0] = "DAY_APPLICATION_value_24"
$data[140,1] = 0.005680248066891133
$data[140,2] = 0.005680248066891133
$data[141,0] = "DAY_APPLICATION_value_9"
$data[141,1] = 0.005378238832675145
$data[141,2] = 0.005378238832675145
$data[142,0] = "DAY_APPLICATION_value_8"
$data[142,1] = 0.005155756922577906
$data[142,2] = -0.005155756922577906
$data[143,0] = "DAY_APPLICATION_value_4"
$data[143,1] = 0.004956650993975522
$data[143,2] = -0.004956650993975522
$data[144,0] = "CATEGORICAL_10_value_HNPAK"
$data[144,1] = 0.004598471426060876
$data[144,2] = -0.004598471426060876
$data[145,0] = "CATEGORICAL_10_value_NFAYV"
$data[145,1] = 0.004598471425962375
$data[145,2] = 0.004598471425962375
$data[146,0] = "CATEGORICAL_0_value_IMFRD"
$data[146,1] = 0.004365475699838027
$data[146,2] = -0.004365475699838027
$data[147,0] = "DOW_APPLICATION_value_1"
$data[147,1] = 0.004051234894245501
$data[147,2] = -0.004051234894245501
$data[148,0] = "DAY_APPLICATION_value_11"
$data[148,1] = 0.003725785357264679
$data[148,2] = -0.003725785357264679
$data[149,0] = "MONTH_APPLICATION_value_1"
$data[149,1] = 0.003656245755848129
$data[149,2] = -0.003656245755848129
$data[150,0] = "DAY_APPLICATION_value_13"
$data[150,1] = 0.003234766395084579
$data[150,2] = -0.003234766395084579
$data[151,0] = "DOW_APPLICATION_value_3"
$data[151,1] = 0.003024219403491201
$data[151,2] = -0.003024219403491201
$data[152,0] = "CATEGORICAL_0_value_UQPEF"
$data[152,1] = 0.002923323083193266
$data[152,2] = -0.002923323083193266
$data[153,0] = "DAY_APPLICATION_value_21"
$data[153,1] = 0.002748099668332766
$data[153,2] = -0.002748099668332766
$data[154,0] = "DOW_APPLICATION_value_2"
$data[154,1] = 0.002587760261459345
$data[154,2] = 0.002587760261459345
$data[155,0] = "DAY_APPLICATION_value_1"
$data[155,1] = 0.002266296641735492
$data[155,2] = -0.002266296641735492
$data[156,0] = "DAY_APPLICATION_value_14"
$data[156,1] = 0.001993389163858855
$data[156,2] = -0.001993389163858855
$data[157,0] = "CATEGORICAL_0_value_FPTCW"
$data[157,1] = 0.001596515774430793
$data[157,2] = 0.001596515774430793
$data[158,0] = "DAY_APPLICATION_value_31"
$data[158,1] = 0.001398359526982454
$data[158,2] = 0.001398359526982454
$data[159,0] = "MONTH_APPLICATION_value_3"
$data[159,1] = 0.0009600351551152781
$data[159,2] = -0.0009600351551152781
$data[160,0] = "MONTH_APPLICATION_value_9"
$data[160,1] = 0.0006891366379682177
$data[160,2] = -0.0006891366379682177
$data[161,0] = "DOW_APPLICATION_value_4"
$data[161,1] = 0.0002945708852986902
$data[161,2] = 0.0002945708852986902
$data[162,0] = "DOW_APPLICATION_value_0"
$data[162,1] = 0.0001101505862857624
$data[162,2] = -0.0001101505862857624
$data[163,0] = "MONTH_APPLICATION_value_5"
$data[163,1] = 0
$data[163,2] = 0

$ws.Range("A2:C165").Value = $data

